# Load_users modal feedback: replace the sample rows of the upload
# template with generic "Usuario Pruebas" / "Prueba Usuario" placeholder
# data instead of real-looking names, and refresh the example
# id/phone numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: Usuario / Pruebas / upruebas / upruebas@yopmail.com / ...
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Usuario"
$ws.Range("B2").Value = "Pruebas"
$ws.Range("C2").Value = "upruebas"
$ws.Range("D2").Value = "upruebas@yopmail.com"
$ws.Range("E2").Value = "Tecnología"
$ws.Range("F2").Value = "Computación"
$ws.Range("H2").Value = "Ingeniería Electrónica y Biomédica"
$ws.Range("I2").Value = "Masculino"
$ws.Range("J2").Value = 20786958
$ws.Range("K2").Value = "Ingeniero"
$ws.Range("L2").Value = "Tecnologías de para el desarrollo"
$ws.Range("N2").Value = 4166542358

# ---------------------------------------------------------------------
# Row 3: Prueba / Usuario / pusuario / pusuario@yopmail.com / ...
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Prueba"
$ws.Range("B3").Value = "Usuario"
$ws.Range("C3").Value = "pusuario"
$ws.Range("D3").Value = "pusuario@yopmail.com"
$ws.Range("E3").Value = "Biociencias"
$ws.Range("F3").Value = "Parasitología"
$ws.Range("H3").Value = "Bioquímica"
$ws.Range("I3").Value = "Femenino"
$ws.Range("J3").Value = 20361852
$ws.Range("K3").Value = "Licenciado "
$ws.Range("L3").Value = "Investigación celular"
$ws.Range("N3").Value = 4127851421

# ---------------------------------------------------------------------
# Hyperlinks: drop the row-3 mailto link entirely, keep only the row-2
# one but pointed at the new sample e-mail address.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:upruebas@yopmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "upruebas@yopmail.com")

# Re-adding the hyperlink resets D2's look to the default "Hyperlink"
# style; restore the original (non-underlined, plain blue) cell format
# by copying it back from D3, which still carries it.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# View: scroll back to the left edge of the sheet and move the active
# selection to K1.
# ---------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$null = $ws.Range("K1").Select()

# ---------------------------------------------------------------------
# Column D width: narrow it slightly (25.98 -> ~21.8 chars).
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 21
